# Applies the "Added ArithUnit and Adder" edit to the Activity Log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Row 7 & 8: fill in the "last 4 digits" column (B) that was previously blank.
$ws.Range("B7").Value = 8414
$ws.Range("B8").Value = 8414

# Row 9: Write LogicUnit.vhd + Set up TBLogicUnit
$ws.Range("B9").Value = 8414
$ws.Range("C9").Value = 43920
$ws.Range("D9").Value = 0.70833333333333337
$ws.Range("E9").Value = 0.73958333333333337
$ws.Range("G9").Value = "Write LogicUnit.vhd + Set up TBLogicUnit"

# Row 10: Write Adder.vhd
$ws.Range("B10").Value = 8414
$ws.Range("C10").Value = 43920
$ws.Range("D10").Value = 0.76041666666666663
$ws.Range("E10").Value = 0.79166666666666663
$ws.Range("G10").Value = "Write Adder.vhd"

# Row 11: Write ArithUnit.vhd
$ws.Range("B11").Value = 8414
$ws.Range("C11").Value = 43920
$ws.Range("D11").Value = 0.79166666666666663
$ws.Range("E11").Value = 0.85416666666666663
$ws.Range("G11").Value = "Write ArithUnit.vhd"

# Update the active selection to reflect where the author left off editing.
$ws.Range("F6").Select()
